$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.314.56"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "3.647.41"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'194.98"
$ws.Range("E5").Value = "  +6.13%  "
$ws.Range("D6").Value = "'578.60"
$ws.Range("D7").Value = "3.639.83"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").Value = "'0.621"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.683"
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("D12").Value = "'57.57"
$ws.Range("E12").Value = "  +7.37%  "
$ws.Range("D13").Value = "'0.0000296"
$ws.Range("E13").Value = "  +16.15%  "
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "4.227.17"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "3.645.96"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("D19").Value = "68.272.23"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").Value = "'405.08"
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").Value = "'12.83"
$ws.Range("E23").Value = "  +24.23%  "
$ws.Range("E24").Value = "  -1.78%  "
$ws.Range("D25").Value = "'86.29"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +3.41%  "
$ws.Range("D27").Value = "'12.68"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").Value = "'3.87"
$ws.Range("E28").Value = "  +7.17%  "
$ws.Range("D29").Value = "'6.09"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  +20.66%  "
$ws.Range("D31").Value = "'9.22"
$ws.Range("E31").Value = "  +2.16%  "
$ws.Range("D32").Value = "'31.81"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("D33").Value = "'692.33"
$ws.Range("E33").Value = "  +17.17%  "
$ws.Range("D34").Value = "'12.27"
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("D36").Value = "'64.77"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("D37").Value = "'42.88"
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("D38").Value = "'0.422"
$ws.Range("E38").Value = "  +11.73%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +8.10%  "
$ws.Range("D41").Value = "'2.86"
$ws.Range("E41").Value = "  +19.51%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "'3.14"
$ws.Range("E42").Value = "  +12.67%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.202.84"
$ws.Range("E43").Value = "  +18.21%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.136"
$ws.Range("E44").Value = "  +2.72%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'2.94"
$ws.Range("E46").Value = "  +28.51%  "
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("E49").Value = "  +8.12%  "
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "'142.89"
$ws.Range("E51").Value = "  +3.92%  "
